$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the Name of the D1 Schottky diode row (E6) to the new part number.
$ws.Range("E6").Value = "BAT54W-HG3-18 "

# The old explicit "General" number style on C4/C8 is no longer needed;
# match the plain (unstyled) format already used elsewhere in the same
# rows (column B carries the default style) so the cells fall back to the
# default style instead of keeping a stray "applyNumberFormat" xf around.
$ws.Range("C4").Style = $ws.Range("B4").Style
$ws.Range("C8").Style = $ws.Range("B8").Style

# Move the active selection to E6, matching where the edit was made.
[void]$ws.Range("E6").Select()
